$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

# Update the title text
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Conclusions & Future Work"

# Append a new paragraph to the content placeholder
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.InsertAfter("`rFuture work could include doing the same analysis on different factor and/or response data, and using daily instead of monthly returns.")
